$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the source data which stores
# prices as literal strings (e.g. "2.348.55") rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.977.92"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "2.351.27"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "545.32"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").Value = "137.09"
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  -8.40%  "
$ws.Range("D9").Value = "2.352.56"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").Value = "5.32"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "0.342"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "24.72"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").Value = "2.776.11"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "60.826.44"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "0.0000160"
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("D18").Value = "2.344.30"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").Value = "10.61"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").Value = "319.75"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "4.12"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").Value = "6.55"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("D25").Value = "63.29"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "8.34"
$ws.Range("E26").Value = "  +8.97%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.467.07"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "7.98"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "502.53"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.38"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0872"
$ws.Range("E32").Value = "  -5.13%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.146"
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.51"
$ws.Range("E35").Value = "  -3.06%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "4.63"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.377"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "18.55"
$ws.Range("E39").Value = "  +3.39%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "5.30"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  +7.29%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "142.34"
$ws.Range("E42").Value = "  +3.82%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "40.54"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "142.72"
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "3.57"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.07"
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0518"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "19.13"
$ws.Range("E48").Value = "  -5.12%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.569"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.0904"
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0220"
$ws.Range("E51").Value = "  -0.88%  "
